# Apply "Target for prev infectious TB" edit to the parameters workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")

# E2: 4 -> 5
$ws.Range("E2").Value = 5

# Row 5 (recent_detection_rate): add distribution/param1/param2
$ws.Range("C5").Value = "uniform"
$ws.Range("D5").Value = 0.1
$ws.Range("E5").Value = 10

# Row 21 (infectiousness_gain_rate): add distribution/param1/param2
$ws.Range("C21").Value = "uniform"
$ws.Range("D21").Value = 0.5
$ws.Range("E21").Value = 5

# Update the sheet view: scroll so row 11 is at the top and select E3
$ws.Activate()
$ws.Range("A11").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("E3").Select() | Out-Null
